$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item(1)
$wsZhCn     = $wb.Worksheets.Item(2)
$wsDeDe     = $wb.Worksheets.Item(3)

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# --- zh-cn sheet: status, handback datetime, error detail ---
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("K2").Value = "2016-08-20 00:53:37"
$wsZhCn.Range("P2").Value = ""

# --- de-de sheet: status, handback datetime, error detail ---
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-08-20 00:53:44"
$wsDeDe.Range("P2").Value = ""

# --- Column width adjustments (best effort, the engine quantizes ColumnWidth) ---
# Overview columns E (5) and F (6) grew to fit the longer status text.
$wsOverview.Columns.Item(5).ColumnWidth = 29.17
$wsOverview.Columns.Item(6).ColumnWidth = 29.17

# zh-cn / de-de column C (3, Status) grew to fit the longer status text,
# column P (16, Error Detail) shrank now that it is empty.
$wsZhCn.Columns.Item(3).ColumnWidth = 29.17
$wsZhCn.Columns.Item(16).ColumnWidth = 12.84

$wsDeDe.Columns.Item(3).ColumnWidth = 29.17
$wsDeDe.Columns.Item(16).ColumnWidth = 12.84
